$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 entries
$ws.Range("A13").Value = 42990
$ws.Range("A13").NumberFormat = "m/d/yy"

$ws.Range("B13").Value = "DiscordiaAgency_Demo_2017_09_12.exe"
$ws.Range("C13").Value = "Entwicklung"
$ws.Range("D13").Value = "Anna Franziska"
$ws.Range("E13").Value = "A* ist drin; Wachen können patrouillieren, zum Gegenstand laufen, Spieler jagen; Schießen fehlt noch"

# Match formatting used by the rest of the sheet for the new row
$ws.Range("A13:E13").VerticalAlignment = -4160

$ws.Range("E13").WrapText = $true

$ws.Rows.Item(13).RowHeight = 60

$ws.Range("C14").Select()
